$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Update the F-column (time_taken) timestamps on the "data" sheet.
# ---------------------------------------------------------------------------
$timestamps = @(
    "2021-10-05 14:33:49.096537",
    "2021-10-05 14:33:49.096549",
    "2021-10-05 14:33:49.096553",
    "2021-10-05 14:33:49.096556",
    "2021-10-05 14:33:49.096559",
    "2021-10-05 14:33:49.096562",
    "2021-10-05 14:33:49.096565",
    "2021-10-05 14:33:49.096568",
    "2021-10-05 14:33:49.096571",
    "2021-10-05 14:33:49.096574",
    "2021-10-05 14:33:49.096577",
    "2021-10-05 14:33:49.096580",
    "2021-10-05 14:33:49.096583",
    "2021-10-05 14:33:49.096586",
    "2021-10-05 14:33:49.096588",
    "2021-10-05 14:33:49.096591",
    "2021-10-05 14:33:49.096594",
    "2021-10-05 14:33:49.096597",
    "2021-10-05 14:33:49.096600",
    "2021-10-05 14:33:49.096602",
    "2021-10-05 14:33:49.096605",
    "2021-10-05 14:33:49.096608",
    "2021-10-05 14:33:49.096611",
    "2021-10-05 14:33:49.096614",
    "2021-10-05 14:33:49.096617",
    "2021-10-05 14:33:49.096620",
    "2021-10-05 14:33:49.096623",
    "2021-10-05 14:33:49.096625",
    "2021-10-05 14:33:49.096628",
    "2021-10-05 14:33:49.096631",
    "2021-10-05 14:33:49.096634",
    "2021-10-05 14:33:49.096636",
    "2021-10-05 14:33:49.096639",
    "2021-10-05 14:33:49.096642",
    "2021-10-05 14:33:49.096645",
    "2021-10-05 14:33:49.096648",
    "2021-10-05 14:33:49.096651",
    "2021-10-05 14:33:49.096654",
    "2021-10-05 14:33:49.096656"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}

# ---------------------------------------------------------------------------
# 2. Add the new "metadata" sheet right after "data".
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# Copy the header style (bold, centered, bordered) from the "data" sheet so
# the style table itself does not need to grow.
$dataSheet.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Epidermolysis bullosa"
$meta.Range("C2").Value = 101
$meta.Range("D2").Value = "1.1"
$meta.Range("D2").NumberFormat = "@"
$meta.Range("E2").Value = "2021-03-11T23:10:09.830174Z"
$meta.Range("F2").Value = "2021-10-05 14:33:49.091913"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/101/?format=json"

# Restore "data" as the active sheet/tab (matches the unchanged activeTab=0).
$dataSheet.Activate()

Write-Host "done"
